$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data table with a new "2021" column (R), mirroring the
# formatting already used by the preceding "2020" column (Q).
$ws.Range("Q4:Q14").Copy() | Out-Null
$ws.Range("R4:R14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 111.17903216128188
$ws.Range("R6").Value = 113.69236134930286
$ws.Range("R7").Value = 114.88854111210361
$ws.Range("R8").Value = 110.91060220352473
$ws.Range("R9").Value = 113.02233875668462
$ws.Range("R10").Value = 110.66816227588356
$ws.Range("R11").Value = 111.40708764208969
$ws.Range("R12").Value = 109.49389157333138
$ws.Range("R13").Value = 110.97185980126036
$ws.Range("R14").Value = 110.008558587758

# Match the author's final selection state.
$ws.Range("T6").Select() | Out-Null
